# Update cryptos list values
# commit: Updated cryptos list on Sat Aug 10 09:48:43 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.855.46"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").Value = "'2.620.20"
$ws.Range("E3").Value = "  -0.04%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'513.98"
$ws.Range("E5").Value = "  +1.10%  "

$ws.Range("D6").Value = "'155.38"
$ws.Range("E6").Value = "  -1.48%  "

$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("E8").Value = "  -0.44%  "

$ws.Range("D9").Value = "'2.635.49"
$ws.Range("E9").Value = "  -0.88%  "

$ws.Range("D10").Value = "'6.81"
$ws.Range("E10").Value = "  +5.32%  "

$ws.Range("E11").Value = "  -0.28%  "

$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("E13").Value = "  +1.81%  "

$ws.Range("D14").Value = "'3.080.00"
$ws.Range("E14").Value = "  +0.12%  "

$ws.Range("D15").Value = "'60.814.55"
$ws.Range("E15").Value = "  +0.48%  "

$ws.Range("D16").Value = "'21.75"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("E17").Value = "  +0.05%  "

$ws.Range("D18").Value = "'2.628.50"
$ws.Range("E18").Value = "  -0.84%  "

$ws.Range("D19").Value = "'4.75"
$ws.Range("E19").Value = "  -0.94%  "

$ws.Range("D20").Value = "'355.96"
$ws.Range("E20").Value = "  +2.75%  "

$ws.Range("E21").Value = "  +1.11%  "

$ws.Range("D22").Value = "'6.20"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("E23").Value = "  +0.28%  "

$ws.Range("D24").Value = "'60.95"
$ws.Range("E24").Value = "  +1.16%  "

$ws.Range("D25").Value = "'0.425"
$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("D26").Value = "'0.167"
$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("E27").Value = "  +0.42%  "

$ws.Range("E28").Value = "  -1.72%  "

$ws.Range("E29").Value = "  -2.91%  "

$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("E31").Value = "  -0.35%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'152.29"
$ws.Range("E32").Value = "  -2.65%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.58"
$ws.Range("E33").Value = "  +0.41%  "

$ws.Range("D34").Value = "'5.93"
$ws.Range("E34").Value = "  +2.53%  "

$ws.Range("D35").Value = "'4.01"
$ws.Range("E35").Value = "  -1.11%  "

$ws.Range("E36").Value = "  -1.20%  "

$ws.Range("D37").Value = "'0.875"
$ws.Range("E37").Value = "  +4.34%  "

$ws.Range("E38").Value = "  +0.14%  "

$ws.Range("E39").Value = "  -0.86%  "

$ws.Range("D40").Value = "'36.38"
$ws.Range("E40").Value = "  +2.76%  "

$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").Value = "'294.25"
$ws.Range("E42").Value = "  -4.61%  "

$ws.Range("E43").Value = "  +0.97%  "

$ws.Range("D44").Value = "'0.626"
$ws.Range("E44").Value = "  -1.73%  "

$ws.Range("D45").Value = "'0.0557"
$ws.Range("E45").Value = "  -3.11%  "

$ws.Range("D46").Value = "'0.996"
$ws.Range("E46").Value = "  +0.50%  "

$ws.Range("D47").Value = "'19.90"
$ws.Range("E47").Value = "  -1.23%  "

$ws.Range("E48").Value = "  +0.48%  "

$ws.Range("E49").Value = "  -0.81%  "

$ws.Range("D50").Value = "'10.30"
$ws.Range("E50").Value = "  +0.20%  "

$ws.Range("D51").Value = "'2.003.38"
$ws.Range("E51").Value = "  -2.49%  "

